$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 14:20"

# Update Canary Islands provinces: Casos totales (B) and Casos activos (C)
# are reported at the archipelago level, so all six island-provinces share
# the same updated totals.
$canaryRows = 23, 24, 25, 26, 27, 28
foreach ($r in $canaryRows) {
    $ws.Cells.Item($r, 2).Value = 1380
    $ws.Cells.Item($r, 3).Value = 62
}

# Update Ceuta: Casos totales (B) and Recuperados (D)
$ceutaRow = 60
$ws.Cells.Item($ceutaRow, 2).Value = 51
$ws.Cells.Item($ceutaRow, 4).Value = 50

# Re-sort the data table (rows 4-64) by "Casos totales" (column B), descending,
# to reflect the updated case counts. Range includes the header row (row 3);
# Header:=xlYes (1) keeps it pinned in place.
$sortRange = $ws.Range("A3:E64")
$sortRange.Sort($ws.Range("B3"), 2, $null, $null, 1, $null, 1, 1)
